$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "C3"   = -11.6323
    "B12"  = 5.063099999999999
    "C14"  = -13.3558
    "C26"  = -11.952
    "C31"  = -13.6088
    "B32"  = 6.745999999999998
    "C35"  = -12.1021
    "B36"  = 9.114700000000001
    "C37"  = -12.8896
    "B38"  = 5.559899999999996
    "C45"  = -13.51489999999999
    "B46"  = 7.013700000000008
    "B54"  = 4.760399999999999
    "B55"  = 5.564
    "C57"  = -14.09869999999999
    "B67"  = 5.544899999999997
    "B69"  = 5.320399999999998
    "B72"  = 5.312700000000002
    "B91"  = 5.388700000000001
    "B99"  = 4.530999999999998
    "C100" = -12.9591
    "C102" = -12.7276
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$wb.Save()
